$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.296878333333333
$ws.Range("H2").Value = 3.890635
$ws.Range("I2").Value = 0.01774073260139904
$ws.Range("J2").Value = 0.02506266560199287
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1011536666666667
$ws.Range("N2").Value = 0.303461
$ws.Range("O2").Value = 0.007629860605400263
$ws.Range("P2").Value = 0.008254451482408482
$ws.Range("Q2").Value = 0.1311839986372222
$ws.Range("R2").Value = 1.180655987735
$ws.Range("S2").Value = 0.0001353593167863547
$ws.Range("T2").Value = 0.0002068785572314781
$ws.Range("G3").Value = 1.296878333333333
$ws.Range("H3").Value = 3.890635
$ws.Range("I3").Value = 0.01774073260139904
$ws.Range("J3").Value = 0.02506266560199287
$ws.Range("O3").Value = 0.7653686681256785
$ws.Range("P3").Value = 0.8280227993585454
$ws.Range("Q3").Value = 13.15936522422278
$ws.Range("R3").Value = 118.434287018005
$ws.Range("S3").Value = 0.01357820088270659
$ws.Range("T3").Value = 0.02075245853114926
$ws.Range("G4").Value = 1.296878333333333
$ws.Range("H4").Value = 3.890635
$ws.Range("I4").Value = 0.01774073260139904
$ws.Range("J4").Value = 0.02506266560199287
$ws.Range("M4").Value = 3.0094955
$ws.Range("N4").Value = 6.018991
$ws.Range("O4").Value = 0.2270014712689213
$ws.Range("P4").Value = 0.1637227491590462
$ws.Range("Q4").Value = 3.902949508214167
$ws.Range("R4").Value = 23.417697049285
$ws.Range("S4").Value = 0.004027172401906099
$ws.Range("T4").Value = 0.004103328513612134
$ws.Range("I5").Value = 0.09943605305674341
$ws.Range("J5").Value = 0.1404751766759988
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1011536666666667
$ws.Range("N5").Value = 0.303461
$ws.Range("O5").Value = 0.007629860605400263
$ws.Range("P5").Value = 0.008254451482408482
$ws.Range("Q5").Value = 0.7352807430093331
$ws.Range("R5").Value = 6.617526687083999
$ws.Range("S5").Value = 0.000758683223974137
$ws.Range("T5").Value = 0.001159545530354792
$ws.Range("I6").Value = 0.09943605305674341
$ws.Range("J6").Value = 0.1404751766759988
$ws.Range("O6").Value = 0.7653686681256785
$ws.Range("P6").Value = 0.8280227993585454
$ws.Range("S6").Value = 0.07610523949171401
$ws.Range("T6").Value = 0.1163166490316468
$ws.Range("I7").Value = 0.09943605305674341
$ws.Range("J7").Value = 0.1404751766759988
$ws.Range("M7").Value = 3.0094955
$ws.Range("N7").Value = 6.018991
$ws.Range("O7").Value = 0.2270014712689213
$ws.Range("P7").Value = 0.1637227491590462
$ws.Range("Q7").Value = 21.875866295734
$ws.Range("R7").Value = 131.255197774404
$ws.Range("S7").Value = 0.02257213034105527
$ws.Range("T7").Value = 0.02299898211399725
$ws.Range("G8").Value = 0.09795233333333332
$ws.Range("H8").Value = 0.293857
$ws.Range("I8").Value = 0.001339945397100812
$ws.Range("J8").Value = 0.001892965987764162
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1011536666666667
$ws.Range("N8").Value = 0.303461
$ws.Range("O8").Value = 0.007629860605400263
$ws.Range("P8").Value = 0.008254451482408482
$ws.Range("Q8").Value = 0.00990823767522222
$ws.Range("R8").Value = 0.08917413907699999
$ws.Range("S8").Value = 0.00001022359659872689
$ws.Range("T8").Value = 0.00001562539590384872
$ws.Range("G9").Value = 0.09795233333333332
$ws.Range("H9").Value = 0.293857
$ws.Range("I9").Value = 0.001339945397100812
$ws.Range("J9").Value = 0.001892965987764162
$ws.Range("O9").Value = 0.7653686681256785
$ws.Range("P9").Value = 0.8280227993585454
$ws.Range("Q9").Value = 0.9939178531767777
$ws.Range("R9").Value = 8.945260678591
$ws.Range("S9").Value = 0.001025552223940182
$ws.Range("T9").Value = 0.001567418996278995
$ws.Range("G10").Value = 0.09795233333333332
$ws.Range("H10").Value = 0.293857
$ws.Range("I10").Value = 0.001339945397100812
$ws.Range("J10").Value = 0.001892965987764162
$ws.Range("M10").Value = 3.0094955
$ws.Range("N10").Value = 6.018991
$ws.Range("O10").Value = 0.2270014712689213
$ws.Range("P10").Value = 0.1637227491590462
$ws.Range("Q10").Value = 0.2947871063811666
$ws.Range("R10").Value = 1.768722638287
$ws.Range("S10").Value = 0.0003041695765619032
$ws.Range("T10").Value = 0.0003099215955813179
$ws.Range("G11").Value = 64.06892400000001
$ws.Range("H11").Value = 128.137848
$ws.Range("I11").Value = 0.8764350668284411
$ws.Range("J11").Value = 0.8254375019458241
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1011536666666667
$ws.Range("N11").Value = 0.303461
$ws.Range("O11").Value = 0.007629860605400263
$ws.Range("P11").Value = 0.008254451482408482
$ws.Range("Q11").Value = 6.480806581988
$ws.Range("R11").Value = 38.88483949192801
$ws.Range("S11").Value = 0.00668707738958567
$ws.Range("T11").Value = 0.006813533811572263
$ws.Range("G12").Value = 64.06892400000001
$ws.Range("H12").Value = 128.137848
$ws.Range("I12").Value = 0.8764350668284411
$ws.Range("J12").Value = 0.8254375019458241
$ws.Range("O12").Value = 0.7653686681256785
$ws.Range("P12").Value = 0.8280227993585454
$ws.Range("Q12").Value = 650.1044460138041
$ws.Range("R12").Value = 3900.626676082825
$ws.Range("S12").Value = 0.6707959397971239
$ws.Range("T12").Value = 0.6834810710567061
$ws.Range("G13").Value = 64.06892400000001
$ws.Range("H13").Value = 128.137848
$ws.Range("I13").Value = 0.8764350668284411
$ws.Range("J13").Value = 0.8254375019458241
$ws.Range("M13").Value = 3.0094955
$ws.Range("N13").Value = 6.018991
$ws.Range("O13").Value = 0.2270014712689213
$ws.Range("P13").Value = 0.1637227491590462
$ws.Range("Q13").Value = 192.815138467842
$ws.Range("R13").Value = 771.2605538713681
$ws.Range("S13").Value = 0.1989520496417314
$ws.Range("T13").Value = 0.1351428970775459
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.3690323333333334
$ws.Range("H14").Value = 1.107097
$ws.Range("I14").Value = 0.005048202116315478
$ws.Range("J14").Value = 0.007131689788420014
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1011536666666667
$ws.Range("N14").Value = 0.303461
$ws.Range("O14").Value = 0.007629860605400263
$ws.Range("P14").Value = 0.008254451482408482
$ws.Range("Q14").Value = 0.03732897363522222
$ws.Range("R14").Value = 0.335960762717
$ws.Range("S14").Value = 0.00003851707845537371
$ws.Range("T14").Value = 0.00005886818734610102
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.3690323333333334
$ws.Range("H15").Value = 1.107097
$ws.Range("I15").Value = 0.005048202116315478
$ws.Range("J15").Value = 0.007131689788420014
$ws.Range("O15").Value = 0.7653686681256785
$ws.Range("P15").Value = 0.8280227993585454
$ws.Range("Q15").Value = 3.744554233856778
$ws.Range("R15").Value = 33.700988104711
$ws.Range("S15").Value = 0.003863735730193609
$ws.Range("T15").Value = 0.005905201742764292
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.3690323333333334
$ws.Range("H16").Value = 1.107097
$ws.Range("I16").Value = 0.005048202116315478
$ws.Range("J16").Value = 0.007131689788420014
$ws.Range("M16").Value = 3.0094955
$ws.Range("N16").Value = 6.018991
$ws.Range("O16").Value = 0.2270014712689213
$ws.Range("P16").Value = 0.1637227491590462
$ws.Range("Q16").Value = 1.110601146521167
$ws.Range("R16").Value = 6.663606879126999
$ws.Range("S16").Value = 0.001145949307666495
$ws.Range("T16").Value = 0.001167619858309621
